# GRAFS example workbook update:
#  - excretion sheet: add a new "Nitrogen Content (%)" column (H) with values
#  - prod sheet: move selection, drop tab-selected state
#  - global sheet: becomes the active/selected tab, new selection
#
# Sheet tab order (0-based): crops, livestock, pop, excretion, prod, global

$wb = $excel.ActiveWorkbook

# --- excretion sheet: new column H "Nitrogen Content (%)" ---
$wsExcretion = $wb.Worksheets.Item("excretion")

$wsExcretion.Range("H1").Value = "Nitrogen Content (%)"
$wsExcretion.Range("H2").Value = 0.5
$wsExcretion.Range("H3").Value = 0.3
$wsExcretion.Range("H4").Value = 0
$wsExcretion.Range("H5").Value = 0.8
$wsExcretion.Range("H6").Value = 0.5
$wsExcretion.Range("H7").Value = 0

# Target OOXML column width is 19.81640625 characters; ColumnWidth=19
# round-trips (via this host's char->stored-width quantization) to the
# closest reachable stored width (19.8333...).
$wsExcretion.Columns.Item(8).ColumnWidth = 19

# Selection ends up on I4 for this sheet (not the active tab though)
$null = $wsExcretion.Range("I4").Select()

# --- prod sheet: selection moves, loses tab focus ---
$wsProd = $wb.Worksheets.Item("prod")
$null = $wsProd.Range("D11").Select()

# --- global sheet: becomes the active tab with a new selection ---
$wsGlobal = $wb.Worksheets.Item("global")
$null = $wsGlobal.Range("B15").Select()
